# Replace the sample "Chetan Parmar" row with a new "Kelvin Kho" row,
# add a new "whatsapp" column (F), and clear the old third data row
# (keeping its formatting) — matches the "added sample card file" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell for the extra column.
$ws.Range("F1").Value = "whatsapp"

# Overwrite row 2 with the new contact's data.
$ws.Range("A2").Value = "Kelvin"
$ws.Range("B2").Value = "Kho"
$ws.Range("C2").Value = "kelvin@techlaju.com"
$ws.Range("D2").Value = "IT Support"
$ws.Range("E2").Value = 60128838318
$ws.Range("F2").Value = 60128838318

# The email cell picks up its own (distinct) style.
$ws.Range("C2").Interior.ColorIndex = -4142

# Row 3's old sample data is gone, but the row keeps its formatting.
$ws.Range("A3:E3").ClearContents()

# Column sizing tweaks for the new layout.
$ws.Columns.Item(3).ColumnWidth = 27.85546875 - 0.8333333333333333
$ws.Columns.Item(4).ColumnWidth = 19 - 0.8333333333333333

# View state: zoomed in, with F2 as the active selection.
$excel.ActiveWindow.Zoom = 205
[void]$ws.Range("F2").Select()
